$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, pushing existing rows 3..21 down to 4..22
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new test data
$ws.Range("A3").Value = "ANC"
$ws.Range("B3").Value = "SEA"
$ws.Range("C3").Value = "GOLDSTREAK"
$ws.Range("D3").Value = "NONSCR"
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 11377
$ws.Range("H3").Value = 11377
$ws.Range("I3").Value = 11377
$ws.Range("J3").Value = "Yes"
$ws.Range("M3").Value = "CAP018_BKG_00001"

# Update current selection to D3, as captured in the workbook after editing
$ws.Range("D3").Select()
